# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '63.360.60'
Set-TextValue 'E2' '  +1.61%  '
Set-TextValue 'D3' '3.168.68'
Set-TextValue 'E3' '  -0.51%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '602.47'
Set-TextValue 'E5' '  +2.79%  '
Set-TextValue 'D6' '135.82'
Set-TextValue 'E6' '  +0.59%  '
Set-TextValue 'E7' '  -0.01%  '
Set-TextValue 'D8' '3.164.62'
Set-TextValue 'E8' '  -0.53%  '
Set-TextValue 'E9' '  +2.10%  '
Set-TextValue 'D10' '0.142'
Set-TextValue 'E10' '  +0.70%  '
Set-TextValue 'D11' '5.35'
Set-TextValue 'E11' '  +2.28%  '
Set-TextValue 'E12' '  +0.97%  '
Set-TextValue 'D13' '0.0000239'
Set-TextValue 'E13' '  +2.08%  '
Set-TextValue 'D14' '34.82'
Set-TextValue 'E14' '  +4.90%  '
Set-TextValue 'D15' '3.692.48'
Set-TextValue 'E15' '  -0.49%  '
Set-TextValue 'E16' '  +1.38%  '
Set-TextValue 'D17' '3.172.77'
Set-TextValue 'E17' '  -0.38%  '
Set-TextValue 'D18' '63.396.41'
Set-TextValue 'E18' '  +1.57%  '
Set-TextValue 'D19' '6.58'
Set-TextValue 'E19' '  +0.08%  '
Set-TextValue 'D20' '461.03'
Set-TextValue 'E20' '  +1.05%  '
Set-TextValue 'D21' '13.97'
Set-TextValue 'E21' '  +0.55%  '
Set-TextValue 'D22' '0.697'
Set-TextValue 'E22' '  -0.85%  '
Set-TextValue 'D23' '7.64'
Set-TextValue 'E23' '  +0.60%  '
Set-TextValue 'E24' '  -0.27%  '
Set-TextValue 'D25' '83.16'
Set-TextValue 'E25' '  +0.79%  '
Set-TextValue 'E26' '  -0.07%  '
Set-TextValue 'E27' '  +0.73%  '
Set-TextValue 'E28' '  -0.06%  '
Set-TextValue 'D29' '2.08'
Set-TextValue 'E29' '  +3.76%  '
Set-TextValue 'B30' 'RenderToken'
Set-TextValue 'C30' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D30' '7.69'
Set-TextValue 'E30' '  -1.49%  '
Set-TextValue 'B31' 'NEARProtocol'
Set-TextValue 'C31' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D31' '6.76'
Set-TextValue 'E31' '  -1.84%  '
Set-TextValue 'D32' '27.11'
Set-TextValue 'E32' '  -0.51%  '
Set-TextValue 'E33' '  -2.42%  '
Set-TextValue 'D34' '2.42'
Set-TextValue 'E34' '  +1.69%  '
Set-TextValue 'E35' '  -1.62%  '
Set-TextValue 'D36' '5.90'
Set-TextValue 'E36' '  +1.87%  '
Set-TextValue 'D37' '0.0₃0732'
Set-TextValue 'E37' '  +6.22%  '
Set-TextValue 'D38' '51.18'
Set-TextValue 'E38' '  +0.03%  '
Set-TextValue 'D39' '0.0390'
Set-TextValue 'E39' '  +0.86%  '
Set-TextValue 'D40' '8.13'
Set-TextValue 'E40' '  +1.62%  '
Set-TextValue 'D41' '0.112'
Set-TextValue 'E41' '  -0.38%  '
Set-TextValue 'D42' '2.62'
Set-TextValue 'E42' '  +0.09%  '
Set-TextValue 'D43' '391.02'
Set-TextValue 'E43' '  -4.68%  '
Set-TextValue 'D44' '2.796.70'
Set-TextValue 'E44' '  -4.94%  '
Set-TextValue 'E45' '  +0.35%  '
Set-TextValue 'D46' '36.01'
Set-TextValue 'E46' '  +0.13%  '
Set-TextValue 'B47' 'USDe'
Set-TextValue 'C47' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D47' '0.999'
Set-TextValue 'E47' '  +0.01%  '
Set-TextValue 'B48' 'Fetch.AI'
Set-TextValue 'C48' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D48' '2.11'
Set-TextValue 'E48' '  -1.67%  '
Set-TextValue 'D49' '125.94'
Set-TextValue 'E49' '  +2.37%  '
Set-TextValue 'D50' '25.07'
Set-TextValue 'E50' '  -1.62%  '
Set-TextValue 'E51' '  +0.97%  '
